# Update CircadiPy cosinor analysis results (sine_10) with re-run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 22.52000000000008
$ws.Range("K2").Value = 47.5261165347294
$ws.Range("L2").Value = "[42.93927393887364, 52.11295913058517]"
$ws.Range("O2").Value = 1.641552918091964
$ws.Range("P2").Value = "[1.540921321580579, 1.7421845146033492]"
$ws.Range("S2").Value = 50.79538107629481
$ws.Range("T2").Value = "[47.7895075122312, 53.80125464035842]"
$ws.Range("W2").Value = 16.63639639639646
$ws.Range("X2").Value = 16.27571571571577
$ws.Range("Y2").Value = 16.99707707707714

# Row 3
$ws.Range("E3").Value = 23.11000000000017
$ws.Range("K3").Value = 48.83109280088215
$ws.Range("L3").Value = "[43.61157481363027, 54.05061078813402]"
$ws.Range("O3").Value = 1.943447707626119
$ws.Range("P3").Value = "[1.8302371615508113, 2.0566582537014257]"
$ws.Range("S3").Value = 55.34334250459163
$ws.Range("T3").Value = "[52.464325744414374, 58.222359264768876]"
$ws.Range("W3").Value = 15.96186186186198
$ws.Range("X3").Value = 15.54546546546558
$ws.Range("Y3").Value = 16.37825825825838
